$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra data rows (3:5), keeping only row 2 of data
$ws.Rows("3:5").Delete()

# Insert two new columns before the old column C ("تاريخ الإنتاج"),
# shifting it to column E
$ws.Columns("C:D").Insert()

# New header cells for the inserted columns
$ws.Range("C1").Value = "المجموعة"
$ws.Range("D1").Value = "الوردية"

# Update the remaining data row
$ws.Range("A2").Value = 185
$ws.Range("B2").Value = 5000
$ws.Range("C2").Value = "الثالثة"
$ws.Range("D2").Value = "صباحية"

# Force the production date to be stored as text, not an auto-converted
# date serial number, then strip the formatting it picked up so the cell
# stays style-less like the rest of the sheet.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2021-03-01"
$ws.Range("E2").ClearFormats()
